$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the "Meta description" paragraph that originally sat right
#    after the title heading.
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# ------------------------------------------------------------------
# 2. Insert a new bold paragraph "Play Genie Jackpots Free | Review of
#    Features & Gameplay" right before the final "Prompt for DALLE"
#    paragraph.
# ------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$dallePara = $d.Paragraphs.Item($lastIndex)

# Create a fresh, empty paragraph right before the DALLE-prompt paragraph.
$null = $dallePara.Range.InsertParagraphBefore()

# The newly created empty paragraph is now paragraph $lastIndex (the
# DALLE-prompt paragraph shifted one slot down).
$newPara = $d.Paragraphs.Item($lastIndex)
$newRange = $d.Range($newPara.Range.Start, $newPara.Range.End)

$newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Genie Jackpots Free | Review of Features &amp; Gameplay</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$null = $newRange.InsertXML($newParaXml)

# ------------------------------------------------------------------
# 3. Replace the DALLE-prompt paragraph's text with the old meta
#    description copy, while keeping the paragraph's italic formatting.
# ------------------------------------------------------------------
$oldPrompt = 'Prompt for DALLE: Create a feature image that brings out the lively and fun nature of the game "Genie Jackpots". The image should be in cartoon style and include a Maya warrior with glasses who looks excited to be playing the game. The warrior should have a big smile on their face and be surrounded by colorful symbols and the genie from the game. The background should be bright and vibrant, with elements of magic and fantasy. Overall, the image should convey a sense of excitement and adventure, depicting the fun experience players can expect from the game.'
$newDescription = 'Find out everything you need to know about Genie Jackpots, including its bonus features and gameplay. Play this simple game for free today!'

$null = $d.Content.Find.Execute($oldPrompt, $true, $false, $false, $false, $false, $true, 1, $false, $newDescription, 2)
